$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$updates = @{
    6  = -3
    7  = -1
    8  = -1
    16 = -2
    17 = -1
    18 = -4
    19 = -2
    22 = -5
    23 = -5
    29 = -5
    39 = -1
    41 = -7
    42 = 4
    43 = -5
    44 = -6
    45 = 5
    49 = 0
    51 = -3
    57 = -2
    58 = 2
    62 = 0
    67 = -4
    69 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
